$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values for rows 1-28 (A:E) per the redone supervision table ---
$ws.Range("A1").Value = 'what'
$ws.Range("B1").Value = 'when'
$ws.Range("C1").Value = 'with'
$ws.Range("D1").Value = 'where'
$ws.Range("E1").Value = 'why'

$ws.Range("A2").Value = 'PhD in Neuroscience  - \textbf{\textit{Summa Cum Laude}}'
$ws.Range("B2").Value = '2015 - 2018'
$ws.Range("C2").Value = 'Postgraduate'
$ws.Range("D2").Value = '\href{https://www.uv.es/}{Universitat de València}, Spain'
$ws.Range("E2").Value = '\href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}. Supervised together with  Alicia Salvador'

$ws.Range("A3").Value = 'Professional Doctorate in Counselling Psychology'
$ws.Range("B3").Value = '2015 - 2018'
$ws.Range("D3").Value = '\href{https://www.uel.ac.uk/}{University of East London}, UK'
$ws.Range("E3").Value = '\href{https://www.researchgate.net/profile/Francisco-Flores-14}{Francisco Javier Flores}. Supervised together with Lisa Chiara Fellin'

$ws.Range("A4").Value = 'Psychological Research Methods (Evolutionary Psychology) MSc - \textbf{\textit{Distinction}}'
$ws.Range("B4").Value = '2013 - 2014'
$ws.Range("D4").Value = '\href{https://www.stir.ac.uk/}{University of Stirling}, UK'
$ws.Range("E4").Value = 'Julia Sanz-Vidania. Supervised together with \href{https://www.scraigroberts.com/}{S Craig Roberts}'

$ws.Range("A5").Value = 'MSc in Psychology - \textbf{\textit{Distinction}}'
$ws.Range("B5").Value = '2019 - 2020'
$ws.Range("D5").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range("E5").Value = 'Adrián Acosta Guerrero. Supervised together with \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'

$ws.Range("A6").Value = 'Biology'
$ws.Range("B6").Value = '2017 - 2018'
$ws.Range("D6").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range("E6").Value = 'Maria Alejandra Abello Mozo  (2017 - 2018)  - \textbf{\textit{Distinction}}'

$ws.Range("A7").Value = 'Music Pedagogy'
$ws.Range("B7").Value = '2017 - 2019'
$ws.Range("D7").Value = '\href{https://www.upn.edu.co/}{Universidad Pedagógica Nacional}, Colombia'
$ws.Range("E7").Value = 'Natalia Elízabeth Moreno Buitrago (2017 ‑ 2019)  - \textbf{\textit{Distinction}}'

$ws.Range("E8").Value = 'Juan Felipe Pérez Ariza (2017 ‑ 2019)  - \textbf{\textit{Distinction}}'

$ws.Range("A9").Value = 'Psychology'
$ws.Range("B9").Value = 'Since 2015'
$ws.Range("C9").Value = 'Undergraduate'
$ws.Range("D9").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range("E9").Value = 'Andrés Castellanos-Chacón (2017 -2018) - \textbf{\textit{Distinction}}. Teaching supervision (2019 - Present)'

$ws.Range("E10").Value = 'Angie Alejandra Lozano Sanjuan (2021 - 2022) - \textbf{\textit{Distinction}}'

$ws.Range("E11").Value = 'Daniela Martínez Franco (2021 - 2022) - \textbf{\textit{Distinction}}'

$ws.Range("E12").Value = 'Mariana Saavedra Botero (2021 - 2022) - \textbf{\textit{Distinction}}'

$ws.Range("E13").Value = 'John Jairo Rubio (2021 - 2022)'

$ws.Range("E14").Value = 'Maria Daniela Martínez Luna (2020 - 2021)  - \textbf{\textit{Distinction}}'

$ws.Range("E15").Value = 'Juan Sebastián Preciado Ruíz (2020 - 2021)  - \textbf{\textit{Distinction}}'

$ws.Range("E16").Value = 'Maria Paula Moreno Rodríguez (2019 - 2021)'

$ws.Range("E17").Value = 'Andrés Felipe Orozco Serrato (2020 - 2021)'

$ws.Range("E18").Value = 'Danny Ferley Gaitan Rodríguez (2019 - 2020)'

$ws.Range("E19").Value = 'Hasbleidy Gamboa Ordoñez (2019 - 2020)'

$ws.Range("E20").Value = 'Paula Andrea Betancourt Velandia  (2018 - 2019)'

$ws.Range("B21").Value = ' '
$ws.Range("E21").Value = 'Ana Sofía Gómez Castelblanco (2018 - 2019)'

$ws.Range("E22").Value = 'Lina María García Hoyos  (2016 - 2017)'

$ws.Range("E23").Value = 'Angie Liliana Pérez Rodríguez  (2016 - 2018)'

$ws.Range("E24").Value = 'Lina María Morales Sánchez (2016 - 2017)'

$ws.Range("E25").Value = 'Laura Milena Estupiñan Aldana  (2016 - 2017)'

$ws.Range("E26").Value = 'Vanesa Díaz Güiza  (2016 - 2018)'

$ws.Range("E27").Value = 'Cindy Paola Moncada Gómez (2016 - 2017)'

$ws.Range("E28").Value = 'Haydn Ricardo Roldán Morales (2015 - 2016)'

# --- Clear the stray space that moved from B19 to B21 ---
$ws.Range("B19").ClearContents()

# --- Add two new blank rows (39, 40) matching the formatting of row 38 ---
$ws.Range("B38:E38").Copy()
$ws.Range("B39:E40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore selection to A7 as in the saved workbook ---
$ws.Range("A7").Select()

